$d = $word.ActiveDocument

# --- Edit 1: simplify the "UI Controls: Text" paragraph -------------------
# Collapse "Draws text using conchars or, in future, the TTF font engine.
# Ignores size for the time being as all text is (8 * 8) * hudscale cvar
# value in size." down to "Draws text using the TTF font engine ".
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute( `
    "conchars or, in future, the TTF font engine. Ignores size for the time being as all text is (8 * 8) * hudscale cvar value in size.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertStart = $findRange.Start
    $findRange.Text = ""
    $newRun = $d.Range($insertStart, $insertStart)
    $newRun.InsertAfter("the TTF font engine ")
}

# --- Edit 2: add a new "UI Networking" bullet after "UI Event Handling" ---
$eventHandlingRange = $d.Content
$eventHandlingRange.Find.ClearFormatting()
$eventHandlingRange.Find.Execute("UI Event Handling", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Resolve via the paragraph Index so we get a "live" Paragraph object
# (one returned straight off a Find range behaves oddly when read back).
$eventHandlingIndex = $eventHandlingRange.Paragraphs.Item(1).Index
$eventHandlingPara = $d.Paragraphs.Item($eventHandlingIndex)
$eventHandlingPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($eventHandlingIndex + 1)
$newPara.Range.Text = "UI Networking"
